$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'22.410.84"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = "'1.566.83"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'1.001"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = "'284.99"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('D7').Value = "'0.3639"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'48.62"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.70%  '
$ws.Range('D9').Value = "'0.3340"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('D10').Value = "'1.126"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.82%  '
$ws.Range('D11').Value = "'0.07412"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').Value = "'20.76"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.88%  '
$ws.Range('D14').Value = "'5.955"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = "'6.904"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').Value = "'1.567.26"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = "'0.00001105"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.76%  '
$ws.Range('D18').Value = "'88.25"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.99%  '
$ws.Range('D19').Value = "'0.06711"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.70%  '
$ws.Range('E20').Value = '  -0.14%  '
$ws.Range('D21').Value = "'6.349"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').Value = "'12.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').Value = "'22.399.32"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').Value = "'2.381"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('D26').Value = "'2.546"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.56%  '
$ws.Range('D27').Value = "'150.33"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.19%  '
$ws.Range('D28').Value = "'19.39"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.73%  '
$ws.Range('D29').Value = "'5.012"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('D30').Value = "'123.90"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').Value = "'1.745.17"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').Value = "'1.062"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.16%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'6.112"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = "'1.998"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('D35').Value = "'9.824"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').Value = "'0.08276"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.49%  '
$ws.Range('D37').Value = "'0.02410"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.12%  '
$ws.Range('D38').Value = "'0.2230"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.73%  '
$ws.Range('D39').Value = "'0.06407"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = "'5.375"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = "'1.288"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.49%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = "'0.6259"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.22%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = "'11.15"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.72%  '
$ws.Range('D44').Value = "'1.000"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = "'13.82"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.01%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = "'0.6070"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.21%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = "'3.759"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('D48').Value = "'2.030"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.52%  '
$ws.Range('D49').Value = "'124.22"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.22%  '
$ws.Range('D50').Value = "'1.218"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.48%  '
$ws.Range('D51').Value = "'0.07203"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.74%  '
